$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65, pushing the existing row 65 down to row 66
$ws.Rows.Item(65).Insert()

# Fill in the new row 65 with the new cherry record (Lapins)
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44578
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100103
$ws.Range("H65").Value = "Frutos de hueso (carozo)"
$ws.Range("I65").Value = 100103001
$ws.Range("J65").Value = "Cereza"
$ws.Range("K65").Value = "Lapins"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 300
$ws.Range("N65").Value = 8500
$ws.Range("O65").Value = 9000
$ws.Range("P65").Value = 8750
$ws.Range("Q65").Value = "`$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia de Curicó"
$ws.Range("S65").Value = 875
$ws.Range("T65").Value = 10
